$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the global "last-not-charged-since" date (column D) for all data rows 2-44
$newD = Get-Date -Year 2025 -Month 12 -Day 12 -Hour 11 -Minute 18 -Second 14
for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 4).Value = $newD
}

# Rows 20-44: refreshed charging-station records (newer terminals / timestamps)
$ws.Cells.Item(20, 1).Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Cells.Item(20, 2).Value = '705号直流'
$ws.Cells.Item(20, 3).Value = (Get-Date -Year 2025 -Month 12 -Day 9 -Hour 12 -Minute 22 -Second 20)

$ws.Cells.Item(21, 1).Value = '长沙特来电飞狐四方坪南区充电站'
$ws.Cells.Item(21, 2).Value = '201号直流'
$ws.Cells.Item(21, 3).Value = (Get-Date -Year 2025 -Month 12 -Day 10 -Hour 5 -Minute 49 -Second 42)

$ws.Cells.Item(22, 1).Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Cells.Item(22, 2).Value = '209号直流'
$ws.Cells.Item(22, 3).Value = (Get-Date -Year 2025 -Month 12 -Day 10 -Hour 14 -Minute 13 -Second 11)

$ws.Cells.Item(23, 1).Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Cells.Item(23, 2).Value = '904号直流'
$ws.Cells.Item(23, 3).Value = (Get-Date -Year 2025 -Month 12 -Day 10 -Hour 15 -Minute 39 -Second 57)

$ws.Cells.Item(24, 1).Value = '长沙特来电飞狐四方坪南区充电站'
$ws.Cells.Item(24, 2).Value = '306号直流'
$ws.Cells.Item(24, 3).Value = (Get-Date -Year 2025 -Month 12 -Day 10 -Hour 22 -Minute 29 -Second 38)

$ws.Cells.Item(25, 1).Value = '长沙特来电飞狐四方坪东区充电站'
$ws.Cells.Item(25, 2).Value = '004A号直流'
$ws.Cells.Item(25, 3).Value = (Get-Date -Year 2025 -Month 12 -Day 11 -Hour 0 -Minute 30 -Second 39)

$ws.Cells.Item(26, 1).Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Cells.Item(26, 2).Value = 'A02号直流'
$ws.Cells.Item(26, 3).Value = (Get-Date -Year 2025 -Month 12 -Day 11 -Hour 4 -Minute 36 -Second 1)

$ws.Cells.Item(27, 1).Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Cells.Item(27, 2).Value = '604号直流'
$ws.Cells.Item(27, 3).Value = (Get-Date -Year 2025 -Month 12 -Day 11 -Hour 9 -Minute 33 -Second 30)

$ws.Cells.Item(28, 1).Value = '长沙特来电飞狐四方坪东区充电站'
$ws.Cells.Item(28, 2).Value = '003B号直流'
$ws.Cells.Item(28, 3).Value = (Get-Date -Year 2025 -Month 12 -Day 11 -Hour 11 -Minute 58 -Second 44)

$ws.Cells.Item(29, 1).Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Cells.Item(29, 2).Value = 'B01号直流'
$ws.Cells.Item(29, 3).Value = (Get-Date -Year 2025 -Month 12 -Day 11 -Hour 13 -Minute 5 -Second 38)

$ws.Cells.Item(30, 1).Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Cells.Item(30, 2).Value = 'B02号直流'
$ws.Cells.Item(30, 3).Value = (Get-Date -Year 2025 -Month 12 -Day 11 -Hour 13 -Minute 10 -Second 30)

$ws.Cells.Item(31, 1).Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Cells.Item(31, 2).Value = '703号直流'
$ws.Cells.Item(31, 3).Value = (Get-Date -Year 2025 -Month 12 -Day 11 -Hour 13 -Minute 14 -Second 6)

$ws.Cells.Item(32, 1).Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Cells.Item(32, 2).Value = '905号直流'
$ws.Cells.Item(32, 3).Value = (Get-Date -Year 2025 -Month 12 -Day 11 -Hour 13 -Minute 28 -Second 13)

$ws.Cells.Item(33, 1).Value = '长沙特来电飞狐四方坪东区充电站'
$ws.Cells.Item(33, 2).Value = '905号直流'
$ws.Cells.Item(33, 3).Value = (Get-Date -Year 2025 -Month 12 -Day 11 -Hour 13 -Minute 59 -Second 11)

$ws.Cells.Item(34, 1).Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Cells.Item(34, 2).Value = '107号直流'
$ws.Cells.Item(34, 3).Value = (Get-Date -Year 2025 -Month 12 -Day 11 -Hour 14 -Minute 0 -Second 9)

$ws.Cells.Item(35, 1).Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Cells.Item(35, 2).Value = '903号直流'
$ws.Cells.Item(35, 3).Value = (Get-Date -Year 2025 -Month 12 -Day 11 -Hour 14 -Minute 31 -Second 25)

$ws.Cells.Item(36, 1).Value = '长沙特来电飞狐四方坪东区充电站'
$ws.Cells.Item(36, 2).Value = '006A号直流'
$ws.Cells.Item(36, 3).Value = (Get-Date -Year 2025 -Month 12 -Day 11 -Hour 14 -Minute 56 -Second 57)

$ws.Cells.Item(37, 1).Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Cells.Item(37, 2).Value = '702号直流'
$ws.Cells.Item(37, 3).Value = (Get-Date -Year 2025 -Month 12 -Day 11 -Hour 15 -Minute 9 -Second 18)

$ws.Cells.Item(38, 1).Value = '长沙特来电飞狐四方坪南区充电站'
$ws.Cells.Item(38, 2).Value = '202号直流'
$ws.Cells.Item(38, 3).Value = (Get-Date -Year 2025 -Month 12 -Day 11 -Hour 15 -Minute 24 -Second 38)

$ws.Cells.Item(39, 1).Value = '长沙特来电飞狐四方坪东区充电站'
$ws.Cells.Item(39, 2).Value = '011B号直流'
$ws.Cells.Item(39, 3).Value = (Get-Date -Year 2025 -Month 12 -Day 11 -Hour 16 -Minute 20 -Second 59)

$ws.Cells.Item(40, 1).Value = '长沙特来电飞狐四方坪南区充电站'
$ws.Cells.Item(40, 2).Value = '403号直流'
$ws.Cells.Item(40, 3).Value = (Get-Date -Year 2025 -Month 12 -Day 11 -Hour 17 -Minute 31 -Second 35)

$ws.Cells.Item(41, 1).Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Cells.Item(41, 2).Value = '505号直流'
$ws.Cells.Item(41, 3).Value = (Get-Date -Year 2025 -Month 12 -Day 11 -Hour 18 -Minute 23 -Second 11)

$ws.Cells.Item(42, 1).Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Cells.Item(42, 2).Value = '403号直流'
$ws.Cells.Item(42, 3).Value = (Get-Date -Year 2025 -Month 12 -Day 11 -Hour 18 -Minute 33 -Second 57)

$ws.Cells.Item(43, 1).Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Cells.Item(43, 2).Value = '503号直流'
$ws.Cells.Item(43, 3).Value = (Get-Date -Year 2025 -Month 12 -Day 11 -Hour 19 -Minute 0 -Second 46)

$ws.Cells.Item(44, 1).Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Cells.Item(44, 2).Value = '805号直流'
$ws.Cells.Item(44, 3).Value = (Get-Date -Year 2025 -Month 12 -Day 11 -Hour 19 -Minute 35 -Second 29)

# Rows 45-54: older records dropped -> blank out A-D (keep formatting)
for ($r = 45; $r -le 54; $r++) {
    $ws.Cells.Item($r, 1).Value = $null
    $ws.Cells.Item($r, 2).Value = $null
    $ws.Cells.Item($r, 3).Value = $null
    $ws.Cells.Item($r, 4).Value = $null
}

# Restore the reported selection
$ws.Range("E13").Select()
